# #7249 - Archived entities should always be read-only - changes after review
#
# Adds 4 new "User Right" rows to the "User Rights" sheet:
#   SORMAS_REST, SORMAS_UI, SORMAS_TO_SORMAS_CLIENT, EXTERNAL_VISITS
# Each row lists, per role column (C:AC), whether that right applies
# ("Yes"/"No"), matching the existing green/red formatted cells used
# throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns C..AC (role columns), in sheet order.
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# Existing cells used purely as formatting templates (kept untouched):
#   A153 -> bold "label" style used for column A
#   C153 -> green "Yes" style
#   D153 -> red "No" style
$labelStyleSrc = "A153"
$yesStyleSrc = "C153"
$noStyleSrc = "D153"

# New rows to append, in order, with their per-column Yes/No values.
$newRows = @(
    @{ Row = 157; Name = "SORMAS_REST"; Values = @("No","No","No","No","Yes","Yes","Yes","No","No","Yes","No","Yes","No","No","No","No","No","No","No","Yes","No","No","No","Yes","Yes","Yes","Yes") },
    @{ Row = 158; Name = "SORMAS_UI"; Values = @("Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","Yes","No","Yes","Yes","No","No","No","No","No") },
    @{ Row = 159; Name = "SORMAS_TO_SORMAS_CLIENT"; Values = @("No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","Yes","No") },
    @{ Row = 160; Name = "EXTERNAL_VISITS"; Values = @("No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","Yes","No","No","No") }
)

foreach ($newRow in $newRows) {
    $r = $newRow.Row
    $name = $newRow.Name
    $values = $newRow.Values

    # Column A: bold "User Right" code, same style as other A-column entries.
    $ws.Range($labelStyleSrc).Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null
    $ws.Range("A" + $r).Value = $name

    # Column B: "Description" - for these rows it just repeats the code.
    $ws.Range("B" + $r).Value = $name

    # Columns C..AC: per-role Yes/No, with matching green/red formatting.
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $values[$i]
        $addr = $col + $r

        if ($val -eq "Yes") {
            $ws.Range($yesStyleSrc).Copy() | Out-Null
        } else {
            $ws.Range($noStyleSrc).Copy() | Out-Null
        }
        $ws.Range($addr).PasteSpecial(-4122) | Out-Null
        $ws.Range($addr).Value = $val
    }
}

$excel.CutCopyMode = 0
